# Remove the "Todo Add BRDF Equations" slide (slide 4) — the materials
# description that was moved out.
$p = $ppt.ActivePresentation
$p.Slides.Item(4).Delete()

# Update the cached date field text (2/3/2017 -> 2/13/2017) on the
# Notes Master and Handout Master date placeholders.
$nm = $p.NotesMaster
$nm.HeadersFooters.DateAndTime.Text = "2/13/2017"

$hm = $p.HandoutMaster
$hm.HeadersFooters.DateAndTime.Text = "2/13/2017"
